$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Wiring up the "debug" pins (p13/p14) and the "debug port" module row ---
# p13 (row 10) and p14 (row 11) are now wired to a debug port -> Usage column (C)
$ws.Range("C10").Value2 = "debug"
$ws.Range("C11").Value2 = "debug"

# The "analog pins" requirement row (row 9) is repurposed to "debug port":
# 2 pins, 1 instance instead of 1 pin, 2 instances.
$ws.Range("E9").Value2 = "debug port"
$ws.Range("F9").Value2 = 2
$ws.Range("G9").Value2 = 1
$ws.Range("H9").Formula = "=F9*G9"

# The "pwm out mosfet" row (row 10 of the requirements table) now needs 0
# pins - it's covered by the debug port wiring above instead.
$ws.Range("F10").Value2 = 0
$ws.Range("G10").Value2 = 0
$ws.Range("H10").Formula = "=F10*G10"

# --- Give the total column (H) real formulas for rows 3-8, matching the
# pattern already used lower down in the table, and pick up the matching
# "total" style used by H7:H10 ---
$ws.Range("H7").Copy()
$ws.Range("H3:H6").PasteSpecial(-4122)

$ws.Range("H3").Formula = "=F3*G3"
$ws.Range("H4").Formula = "=F4*G4"
$ws.Range("H5").Formula = "=F5*G5"
$ws.Range("H6").Formula = "=F6*G6"
$ws.Range("H7").Formula = "=F7*G7"
$ws.Range("H8").Formula = "=F8*G8"

# --- Move the active selection to where the edits were made ---
$ws.Range("E11").Select() | Out-Null
